$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot with the latest scraped values.
# Cells are forced to Text format before assignment (and restored to the default
# "Normal" style afterwards) so that values such as "63.994.64" or "0.0000255"
# are stored as literal text instead of being reinterpreted as numbers by Excel.

# Row 2: update D2, E2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.994.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("E2").Style = "Normal"

# Row 3: update D3, E3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.286.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E3").Style = "Normal"

# Row 4: update E4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E4").Style = "Normal"

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("E5").Style = "Normal"

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.46%  '
$ws.Range("E6").Style = "Normal"

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.583'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.56%  '
$ws.Range("E7").Style = "Normal"

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.277.60'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("E8").Style = "Normal"

# Row 9: update E9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E9").Style = "Normal"

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.599'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("E10").Style = "Normal"

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -12.82%  '
$ws.Range("E11").Style = "Normal"

# Row 12: update E12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("E12").Style = "Normal"

# Row 13: update D13, E13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("E13").Style = "Normal"

# Row 14: update D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.27%  '
$ws.Range("E14").Style = "Normal"

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.832.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("E15").Style = "Normal"

# Row 16: update B16, C16, D16, E16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TRON'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.117'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("E16").Style = "Normal"

# Row 17: update B17, C17, D17, E17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.293.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("E17").Style = "Normal"

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.809.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("E18").Style = "Normal"

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("E19").Style = "Normal"

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("E20").Style = "Normal"

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.946'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("E21").Style = "Normal"

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '377.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E22").Style = "Normal"

# Row 23: update E23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.27%  '
$ws.Range("E23").Style = "Normal"

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("E24").Style = "Normal"

# Row 25: update B25, C25, D25, E25
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.37%  '
$ws.Range("E25").Style = "Normal"

# Row 26: update B26, C26, D26, E26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("E26").Style = "Normal"

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E27").Style = "Normal"

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("E28").Style = "Normal"

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.60%  '
$ws.Range("E29").Style = "Normal"

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.90%  '
$ws.Range("E30").Style = "Normal"

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("E31").Style = "Normal"

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '626.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.22%  '
$ws.Range("E32").Style = "Normal"

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("E33").Style = "Normal"

# Row 34: update E34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.28%  '
$ws.Range("E34").Style = "Normal"

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.104'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("E35").Style = "Normal"

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '56.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.09%  '
$ws.Range("E36").Style = "Normal"

# Row 37: update E37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E37").Style = "Normal"

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.77%  '
$ws.Range("E38").Style = "Normal"

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.374'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.01%  '
$ws.Range("E39").Style = "Normal"

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0742'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.99%  '
$ws.Range("E40").Style = "Normal"

# Row 41: update E41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E41").Style = "Normal"

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +11.04%  '
$ws.Range("E42").Style = "Normal"

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.124'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E43").Style = "Normal"

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.08%  '
$ws.Range("E44").Style = "Normal"

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.889.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E45").Style = "Normal"

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("E46").Style = "Normal"

# Row 47: update E47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.17%  '
$ws.Range("E47").Style = "Normal"

# Row 48: update E48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("E48").Style = "Normal"

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.41%  '
$ws.Range("E49").Style = "Normal"

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("E50").Style = "Normal"

# Row 51: update E51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.29%  '
$ws.Range("E51").Style = "Normal"
